# Auto-generated: apply scheduled-runner price/profit refresh to Odin_Profits sheets.
# Each sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) gets updated market-price
# derived columns (H..N) for a handful of leve rows.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 1739
$ws.Range("I31").Value = 1739
$ws.Range("K31").Value = 5217
$ws.Range("M31").Value = -4987
$ws.Range("H62").Value = 23820046
$ws.Range("I62").Value = 55573884
$ws.Range("J62").Value = 4667.5
$ws.Range("K62").Value = 55573884
$ws.Range("L62").Value = 4667.5
$ws.Range("M62").Value = -55573260
$ws.Range("N62").Value = -5915.5
$ws.Range("H65").Value = 23820046
$ws.Range("I65").Value = 55573884
$ws.Range("J65").Value = 4667.5
$ws.Range("K65").Value = 277869420
$ws.Range("L65").Value = 23337.5
$ws.Range("M65").Value = -277866300
$ws.Range("N65").Value = -29577.5
$ws.Range("H70").Value = 14600.75
$ws.Range("J70").Value = 14600.75
$ws.Range("L70").Value = 43802.25
$ws.Range("N70").Value = -44342.25
$ws.Range("H73").Value = 14600.75
$ws.Range("J73").Value = 14600.75
$ws.Range("L73").Value = 43802.25
$ws.Range("N73").Value = -45674.25
$ws.Range("H74").Value = 7944176.5
$ws.Range("I74").Value = 12992970
$ws.Range("K74").Value = 12992970
$ws.Range("M74").Value = -12992034
$ws.Range("H77").Value = 7944176.5
$ws.Range("I77").Value = 12992970
$ws.Range("K77").Value = 64964850
$ws.Range("M77").Value = -64960170
$ws.Range("H100").Value = 3742.3
$ws.Range("I100").Value = 3324.7778
$ws.Range("K100").Value = 3324.7778
$ws.Range("M100").Value = -2783.7778
$ws.Range("H132").Value = 338489.44
$ws.Range("I132").Value = 383560.22
$ws.Range("K132").Value = 1150680.66
$ws.Range("M132").Value = -1148150.66
$ws.Range("H137").Value = 6737.357
$ws.Range("I137").Value = 7447.375
$ws.Range("K137").Value = 22342.125
$ws.Range("M137").Value = -19792.125

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1391.3334
$ws.Range("I45").Value = 1071.625
$ws.Range("J45").Value = 1756.7142
$ws.Range("K45").Value = 1071.625
$ws.Range("L45").Value = 1756.7142
$ws.Range("M45").Value = -694.625
$ws.Range("N45").Value = -2510.7142
$ws.Range("H110").Value = 4800.533
$ws.Range("I110").Value = 2018.5
$ws.Range("J110").Value = 6655.222
$ws.Range("K110").Value = 2018.5
$ws.Range("L110").Value = 6655.222
$ws.Range("M110").Value = 26.5
$ws.Range("N110").Value = -10745.222
$ws.Range("H132").Value = 617185.9399999999
$ws.Range("I132").Value = 744184.9
$ws.Range("J132").Value = 88023.586
$ws.Range("K132").Value = 2232554.7
$ws.Range("L132").Value = 264070.758
$ws.Range("M132").Value = -2230024.7
$ws.Range("N132").Value = -269130.758
$ws.Range("H133").Value = 197983
$ws.Range("J133").Value = 197983
$ws.Range("L133").Value = 197983
$ws.Range("N133").Value = -203043

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 20849998
$ws.Range("J80").Value = 27799120
$ws.Range("L80").Value = 27799120
$ws.Range("N80").Value = -27801116
$ws.Range("H83").Value = 20849998
$ws.Range("J83").Value = 27799120
$ws.Range("L83").Value = 138995600
$ws.Range("N83").Value = -139005584
$ws.Range("H99").Value = 7683.4775
$ws.Range("I99").Value = 7023.48
$ws.Range("J99").Value = 9624.647000000001
$ws.Range("K99").Value = 7023.48
$ws.Range("L99").Value = 9624.647000000001
$ws.Range("M99").Value = -5525.48
$ws.Range("N99").Value = -12620.647

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 25005318
$ws.Range("I16").Value = 33336972
$ws.Range("J16").Value = 10362.4
$ws.Range("K16").Value = 33336972
$ws.Range("L16").Value = 10362.4
$ws.Range("M16").Value = -33336685
$ws.Range("N16").Value = -10936.4
$ws.Range("H113").Value = 25005318
$ws.Range("I113").Value = 33336972
$ws.Range("J113").Value = 10362.4
$ws.Range("K113").Value = 33336972
$ws.Range("L113").Value = 10362.4
$ws.Range("M113").Value = -33334802
$ws.Range("N113").Value = -14702.4
$ws.Range("H121").Value = 107156.25
$ws.Range("J121").Value = 107156.25
$ws.Range("L121").Value = 107156.25
$ws.Range("N121").Value = -109776.25
$ws.Range("H122").Value = 6977.5713
$ws.Range("I122").Value = 2067.2
$ws.Range("K122").Value = 6201.599999999999
$ws.Range("M122").Value = -3751.599999999999
$ws.Range("H132").Value = 7561.75
$ws.Range("I132").Value = 3439.3333
$ws.Range("K132").Value = 10317.9999
$ws.Range("M132").Value = -7787.999899999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 108418.45
$ws.Range("J37").Value = 108418.45
$ws.Range("L37").Value = 325255.35
$ws.Range("N37").Value = -325479.35
$ws.Range("H107").Value = 3746
$ws.Range("I107").Value = 957.6667
$ws.Range("J107").Value = 4238.0586
$ws.Range("K107").Value = 2873.0001
$ws.Range("L107").Value = 12714.1758
$ws.Range("M107").Value = -953.0001000000002
$ws.Range("N107").Value = -16554.1758
$ws.Range("H122").Value = 4378.0386
$ws.Range("I122").Value = 985.8333
$ws.Range("J122").Value = 5395.7
$ws.Range("K122").Value = 8872.4997
$ws.Range("L122").Value = 48561.3
$ws.Range("M122").Value = -6422.4997
$ws.Range("N122").Value = -53461.3

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5418.857
$ws.Range("I70").Value = 4655.5835
$ws.Range("K70").Value = 4655.5835
$ws.Range("M70").Value = -4385.5835
$ws.Range("H73").Value = 5418.857
$ws.Range("I73").Value = 4655.5835
$ws.Range("K73").Value = 4655.5835
$ws.Range("M73").Value = -3719.5835
$ws.Range("H113").Value = 6146.048
$ws.Range("I113").Value = 1703.3334
$ws.Range("K113").Value = 1703.3334
$ws.Range("M113").Value = 466.6666
$ws.Range("H132").Value = 4865.896
$ws.Range("J132").Value = 4217.273
$ws.Range("L132").Value = 12651.819
$ws.Range("N132").Value = -17711.819
$ws.Range("H137").Value = 58000
$ws.Range("J137").Value = 58000
$ws.Range("L137").Value = 58000
$ws.Range("N137").Value = -68200

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4973.2964
$ws.Range("I61").Value = 3981.16
$ws.Range("K61").Value = 3981.16
$ws.Range("M61").Value = -3779.16
$ws.Range("H113").Value = 4973.2964
$ws.Range("I113").Value = 3981.16
$ws.Range("K113").Value = 3981.16
$ws.Range("M113").Value = -1811.16
$ws.Range("H122").Value = 3333249.2
$ws.Range("I122").Value = 19961996
$ws.Range("K122").Value = 59885988
$ws.Range("M122").Value = -59883538

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2085.8
$ws.Range("J107").Value = 3119.25
$ws.Range("L107").Value = 9357.75
$ws.Range("N107").Value = -13197.75
$ws.Range("H136").Value = 18532420
$ws.Range("I136").Value = 33348714
$ws.Range("J136").Value = 12053.583
$ws.Range("K136").Value = 100046142
$ws.Range("L136").Value = 36160.749
$ws.Range("M136").Value = -100043592
$ws.Range("N136").Value = -41260.749
